$wb = $excel.ActiveWorkbook

# Sheet "展览" - update F2 (想去人数) and F4 (想去人数)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 279
$ws1.Range("F4").Value = 1059

# Sheet "全部类型" - same two rows mirrored
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 279
$ws4.Range("F4").Value = 1059
